$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string date labels for the newly appended rows (column A).
$dates = @{
    124 = "02 06 2020"
    125 = "03 06 2020"
    126 = "04 06 2020"
    127 = "05 06 2020"
    128 = "06 06 2020"
    129 = "07 06 2020"
    130 = "08 06 2020"
}

foreach ($r in $dates.Keys) {
    $ws.Range("A$r").Value = $dates[$r]
}

# Numeric data for rows 121-125 (B:BE, skipping the columns that stay blank
# for every row: E, N, AC, AR, AZ - matching the rest of the sheet).
$data = @{
    "B121" = 7.9793633664459
    "C121" = 28.192103022013
    "D121" = 2.7622782844218
    "F121" = 19.26869413428
    "G121" = 9.755931811637099
    "H121" = 7.6944489497781
    "I121" = 14.649260052593
    "J121" = 21.165382416527
    "K121" = 0.73035327311622
    "L121" = 10.209984259055
    "M121" = 19.900484662777
    "O121" = 19.65571081483
    "P121" = 18.531937544448
    "Q121" = 30.083736418867
    "R121" = 18.54789461307
    "S121" = 11.54622864115
    "T121" = 29.274399833841
    "U121" = 17.047242217425
    "V121" = 15.635438323867
    "W121" = 24.235206158101
    "X121" = 17.342062953384
    "Y121" = 9.809186191754501
    "Z121" = 14.375812879596
    "AA121" = 18.952137270384
    "AB121" = 23.423678577915
    "AD121" = 6.5675985118081
    "AE121" = 12.781980498621
    "AF121" = 19.566816146605
    "AG121" = 52.689514679252
    "AH121" = 20.618120962682
    "AI121" = 27.866317107946
    "AJ121" = 15.893127791825
    "AK121" = 4.1684779396163
    "AL121" = 13.859553491327
    "AM121" = 13.873918559562
    "AN121" = 13.187131220059
    "AO121" = 14.08609123705
    "AP121" = 0.49211577166877
    "AQ121" = 17.138834676148
    "AS121" = 26.570557849293
    "AT121" = 12.107879904044
    "AU121" = 28.531612161365
    "AV121" = 28.695993562957
    "AW121" = 12.265849863824
    "AX121" = 23.257677556078
    "AY121" = 21.009251961086
    "BA121" = 29.263906907568
    "BB121" = 11.131767360069
    "BC121" = 13.514553908772
    "BD121" = 12.550292575792
    "BE121" = 5.0904742527433
    "B122" = 3.65004346
    "C122" = 27.8697077
    "D122" = 1.51317463
    "F122" = 19.76359971
    "G122" = 8.265969849999999
    "H122" = 4.36380014
    "I122" = 9.710641839999999
    "J122" = 12.44229112
    "K122" = 0
    "L122" = 7.12809171
    "M122" = 16.25072625
    "O122" = 15.67049883
    "P122" = 13.48300133
    "Q122" = 23.84991925
    "R122" = 16.95679145
    "S122" = 7.78752718
    "T122" = 22.99940487
    "U122" = 11.65694783
    "V122" = 23.884027
    "W122" = 17.92131912
    "X122" = 17.99011484
    "Y122" = 6.43688664
    "Z122" = 10.0878038
    "AA122" = 13.81895744
    "AB122" = 23.44622886
    "AD122" = 4.93398457
    "AE122" = 8.63148221
    "AF122" = 15.96056651
    "AG122" = 38.10297313
    "AH122" = 16.05116721
    "AI122" = 17.10936954
    "AJ122" = 10.70024875
    "AK122" = 0.42056642
    "AL122" = 9.47561163
    "AM122" = 11.70473825
    "AN122" = 8.99938935
    "AO122" = 10.46861077
    "AP122" = 0
    "AQ122" = 12.76758938
    "AS122" = 20.86167863
    "AT122" = 9.00301595
    "AU122" = 23.43453216
    "AV122" = 25.10731323
    "AW122" = 10.96353683
    "AX122" = 26.49324001
    "AY122" = 19.03021295
    "BA122" = 23.50196581
    "BB122" = 11.06733191
    "BC122" = 9.602136740000001
    "BD122" = 19.80941516
    "BE122" = 2.58132174
    "B123" = 0.33261593
    "C123" = 26.72782824
    "D123" = 0.51673582
    "F123" = 16.21266104
    "G123" = 7.11528969
    "H123" = 1.73307085
    "I123" = 11.27197903
    "J123" = 5.43032054
    "K123" = 0
    "L123" = 4.48569053
    "M123" = 14.40723369
    "O123" = 11.84629811
    "P123" = 22.72056603
    "Q123" = 17.91722489
    "R123" = 17.28638513
    "S123" = 11.12361928
    "T123" = 27.70669178
    "U123" = 10.53132546
    "V123" = 19.28559141
    "W123" = 18.32209622
    "X123" = 12.65911024
    "Y123" = 3.58587684
    "Z123" = 10.14289731
    "AA123" = 9.34993148
    "AB123" = 22.51639154
    "AD123" = 3.44367724
    "AE123" = 5.14387324
    "AF123" = 17.1895849
    "AG123" = 25.20269832
    "AH123" = 24.99290531
    "AI123" = 8.36831941
    "AJ123" = 9.631990699999999
    "AK123" = 0
    "AL123" = 17.798114
    "AM123" = 10.29631048
    "AN123" = 7.59109927
    "AO123" = 7.31571323
    "AP123" = 0
    "AQ123" = 12.10834695
    "AS123" = 15.62441736
    "AT123" = 6.25393297
    "AU123" = 36.62207215
    "AV123" = 19.33846922
    "AW123" = 9.54122564
    "AX123" = 20.47635393
    "AY123" = 14.24071164
    "BA123" = 18.0220487
    "BB123" = 9.33711001
    "BC123" = 10.52602255
    "BD123" = 15.08302874
    "BE123" = 0.57841406
    "B124" = 0
    "C124" = 19.54634249
    "D124" = 0
    "F124" = 11.3761136
    "G124" = 6.39803663
    "H124" = 0
    "I124" = 10.13823245
    "J124" = 38.41940162
    "K124" = 0
    "L124" = 4.43171955
    "M124" = 12.29190437
    "O124" = 8.348088860000001
    "P124" = 17.39827087
    "Q124" = 12.514154
    "R124" = 16.14245443
    "S124" = 7.78096814
    "T124" = 21.20781956
    "U124" = 6.43529045
    "V124" = 28.18993689
    "W124" = 12.8850658
    "X124" = 13.33761491
    "Y124" = 1.29674604
    "Z124" = 9.222967560000001
    "AA124" = 7.6618515
    "AB124" = 21.49262175
    "AD124" = 2.13668045
    "AE124" = 2.31855851
    "AF124" = 16.58513926
    "AG124" = 14.32949972
    "AH124" = 19.72882658
    "AI124" = 1.61756365
    "AJ124" = 8.97017224
    "AK124" = 0
    "AL124" = 17.40202088
    "AM124" = 8.85165497
    "AN124" = 4.46281503
    "AO124" = 4.64280044
    "AP124" = 0.97702061
    "AQ124" = 11.13261749
    "AS124" = 10.94501316
    "AT124" = 6.66947289
    "AU124" = 29.58568074
    "AV124" = 16.01379889
    "AW124" = 8.11083241
    "AX124" = 14.96960451
    "AY124" = 11.54717713
    "BA124" = 53.59802094
    "BB124" = 6.41102402
    "BC124" = 9.42868908
    "BD124" = 20.80958069
    "BE124" = 0
    "B125" = 0
    "C125" = 27.53724147
    "D125" = 0
    "F125" = 10.63284612
    "G125" = 5.40077105
    "H125" = 0
    "I125" = 11.98111968
    "J125" = 53.92679791
    "K125" = 0
    "L125" = 6.09252161
    "M125" = 10.41332465
    "O125" = 5.29769845
    "P125" = 12.60070428
    "Q125" = 7.81239225
    "R125" = 15.29151498
    "S125" = 4.95017637
    "T125" = 20.31850067
    "U125" = 9.7244007
    "V125" = 22.54734873
    "W125" = 9.6731011
    "X125" = 12.22125457
    "Y125" = 0
    "Z125" = 8.43037247
    "AA125" = 8.458074119999999
    "AB125" = 20.13981098
    "AD125" = 6.74881888
    "AE125" = 0.13166108
    "AF125" = 13.13324635
    "AG125" = 44.56530818
    "AH125" = 14.8212818
    "AI125" = 0
    "AJ125" = 6.48714289
    "AK125" = 0
    "AL125" = 16.77220248
    "AM125" = 8.26309198
    "AN125" = 8.04869225
    "AO125" = 2.44921995
    "AP125" = 0.41727147
    "AQ125" = 10.23971034
    "AS125" = 6.89401319
    "AT125" = 6.98211387
    "AU125" = 22.79200537
    "AV125" = 14.77535793
    "AW125" = 7.71545812
    "AX125" = 28.85687188
    "AY125" = 12.97377111
    "BA125" = 63.05946682
    "BB125" = 9.713180899999999
    "BC125" = 8.5274167
    "BD125" = 15.90845047
    "BE125" = 0
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
